# Auto-generated script applying scheduled market-data refresh to all 8 Leve sheets
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H4").Value = 7143206.5
$ws.Range("I4").Value = 10000300
$ws.Range("J4").Value = 473.25
$ws.Range("K4").Value = 10000300
$ws.Range("L4").Value = 473.25
$ws.Range("M4").Value = -10000186
$ws.Range("N4").Value = -701.25

$ws.Range("H8").Value = 662
$ws.Range("I8").Value = 263.8
$ws.Range("J8").Value = 1325.6666
$ws.Range("K8").Value = 791.4000000000001
$ws.Range("L8").Value = 3976.9998
$ws.Range("M8").Value = -652.4000000000001
$ws.Range("N8").Value = -4254.9998

$ws.Range("H10").Value = 23134.166
$ws.Range("I10").Value = 12000
$ws.Range("J10").Value = 25361
$ws.Range("K10").Value = 12000
$ws.Range("L10").Value = 25361
$ws.Range("M10").Value = -11707
$ws.Range("N10").Value = -25947

$ws.Range("H19").Value = 449.4138
$ws.Range("I19").Value = 417.84616
$ws.Range("J19").Value = 475.0625
$ws.Range("K19").Value = 417.84616
$ws.Range("L19").Value = 475.0625
$ws.Range("M19").Value = -242.84616
$ws.Range("N19").Value = -825.0625

$ws.Range("H80").Value = 3675.5667
$ws.Range("I80").Value = 359.94446
$ws.Range("J80").Value = 8649
$ws.Range("K80").Value = 1079.83338
$ws.Range("L80").Value = 25947
$ws.Range("M80").Value = -81.83338000000003
$ws.Range("N80").Value = -27943

$ws.Range("H83").Value = 3675.5667
$ws.Range("I83").Value = 359.94446
$ws.Range("J83").Value = 8649
$ws.Range("K83").Value = 3239.50014
$ws.Range("L83").Value = 77841
$ws.Range("M83").Value = 1752.49986
$ws.Range("N83").Value = -87825

$ws.Range("H105").Value = 29444.445
$ws.Range("J105").Value = 29444.445
$ws.Range("L105").Value = 29444.445
$ws.Range("N105").Value = -36432.445

$ws.Range("H138").Value = 2730.405
$ws.Range("I138").Value = 1366.7297
$ws.Range("J138").Value = 3931.738
$ws.Range("K138").Value = 4100.189100000001
$ws.Range("L138").Value = 11795.214
$ws.Range("M138").Value = 1039.810899999999
$ws.Range("N138").Value = -22075.214

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H3").Value = 8499
$ws.Range("I3").Value = 8499
$ws.Range("J3").Value = 0
$ws.Range("K3").Value = 8499
$ws.Range("L3").Value = 0
$ws.Range("M3").Value = -8384
$ws.Range("N3").Value = $null

$ws.Range("H32").Value = 5101.08
$ws.Range("I32").Value = 4362.553
$ws.Range("J32").Value = 16671.334
$ws.Range("K32").Value = 4362.553
$ws.Range("L32").Value = 16671.334
$ws.Range("M32").Value = -4075.553
$ws.Range("N32").Value = -17245.334

$ws.Range("H74").Value = 9616838
$ws.Range("I74").Value = 1234.6897
$ws.Range("J74").Value = 21740860
$ws.Range("K74").Value = 1234.6897
$ws.Range("L74").Value = 21740860
$ws.Range("M74").Value = -360.6896999999999
$ws.Range("N74").Value = -21742608

$ws.Range("H77").Value = 9616838
$ws.Range("I77").Value = 1234.6897
$ws.Range("J77").Value = 21740860
$ws.Range("K77").Value = 6173.4485
$ws.Range("L77").Value = 108704300
$ws.Range("M77").Value = -1805.4485
$ws.Range("N77").Value = -108713036

$ws.Range("H132").Value = 1411019.2
$ws.Range("I132").Value = 1781.7551
$ws.Range("K132").Value = 5345.2653
$ws.Range("M132").Value = -2815.2653

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H132").Value = 51706.668
$ws.Range("J132").Value = 51706.668
$ws.Range("L132").Value = 51706.668
$ws.Range("N132").Value = -61826.668

$ws.Range("H133").Value = 43232.5
$ws.Range("J133").Value = 43694.285
$ws.Range("L133").Value = 43694.285
$ws.Range("N133").Value = -53814.285

$ws.Range("H134").Value = 22673.236
$ws.Range("I134").Value = 4900.0454
$ws.Range("J134").Value = 93766
$ws.Range("K134").Value = 14700.1362
$ws.Range("L134").Value = 281298
$ws.Range("M134").Value = -12165.1362
$ws.Range("N134").Value = -286368

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H58").Value = 6311145.5
$ws.Range("I58").Value = 7938002.5
$ws.Range("J58").Value = 1430575.6
$ws.Range("K58").Value = 7938002.5
$ws.Range("L58").Value = 1430575.6
$ws.Range("M58").Value = -7937799.5
$ws.Range("N58").Value = -1430981.6

$ws.Range("H94").Value = 2578.4546
$ws.Range("I94").Value = 2321
$ws.Range("J94").Value = 2675
$ws.Range("K94").Value = 2321
$ws.Range("L94").Value = 2675
$ws.Range("M94").Value = -1870
$ws.Range("N94").Value = -3577

$ws.Range("H134").Value = 7632826
$ws.Range("I134").Value = 10755748
$ws.Range("K134").Value = 32267244
$ws.Range("M134").Value = -32264709

$ws.Range("H136").Value = 6311145.5
$ws.Range("I136").Value = 7938002.5
$ws.Range("J136").Value = 1430575.6
$ws.Range("K136").Value = 23814007.5
$ws.Range("L136").Value = 4291726.800000001
$ws.Range("M136").Value = -23811457.5
$ws.Range("N136").Value = -4296826.800000001

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H11").Value = 72.8
$ws.Range("I11").Value = 72.8
$ws.Range("K11").Value = 218.4
$ws.Range("M11").Value = -78.39999999999998

$ws.Range("H92").Value = 832.0625
$ws.Range("J92").Value = 836.3570999999999
$ws.Range("L92").Value = 2509.0713
$ws.Range("N92").Value = -5005.0713

$ws.Range("H131").Value = 2440235.2
$ws.Range("J131").Value = 1573.9286
$ws.Range("L131").Value = 4721.7858
$ws.Range("N131").Value = -14801.7858

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H5").Value = 8500
$ws.Range("J5").Value = 8500
$ws.Range("L5").Value = 8500
$ws.Range("N5").Value = -8724

$ws.Range("H12").Value = 6800000
$ws.Range("I12").Value = 6800000
$ws.Range("J12").Value = 0
$ws.Range("K12").Value = 6800000
$ws.Range("L12").Value = 0
$ws.Range("M12").Value = -6799860
$ws.Range("N12").Value = $null

$ws.Range("H22").Value = 0
$ws.Range("I22").Value = 0
$ws.Range("K22").Value = 0
$ws.Range("M22").Value = $null

$ws.Range("H33").Value = 0
$ws.Range("J33").Value = 0
$ws.Range("L33").Value = 0
$ws.Range("N33").Value = $null

$ws.Range("H36").Value = 1240
$ws.Range("I36").Value = 875
$ws.Range("J36").Value = 1483.3334
$ws.Range("K36").Value = 875
$ws.Range("L36").Value = 1483.3334
$ws.Range("M36").Value = -390
$ws.Range("N36").Value = -2453.3334

$ws.Range("H40").Value = 8500
$ws.Range("J40").Value = 8500
$ws.Range("L40").Value = 8500
$ws.Range("N40").Value = -8802

$ws.Range("H43").Value = 3325
$ws.Range("I43").Value = 1600
$ws.Range("K43").Value = 1600
$ws.Range("M43").Value = -1449

$ws.Range("H44").Value = 100000
$ws.Range("J44").Value = 100000
$ws.Range("L44").Value = 100000
$ws.Range("N44").Value = -101192

$ws.Range("H49").Value = 18000
$ws.Range("J49").Value = 18000
$ws.Range("L49").Value = 18000
$ws.Range("N49").Value = -18368

$ws.Range("H53").Value = 20000
$ws.Range("I53").Value = 0
$ws.Range("K53").Value = 0
$ws.Range("M53").Value = $null

$ws.Range("H55").Value = 4001
$ws.Range("I55").Value = 0
$ws.Range("J55").Value = 4001
$ws.Range("K55").Value = 0
$ws.Range("L55").Value = 4001
$ws.Range("M55").Value = $null
$ws.Range("N55").Value = -4655

$ws.Range("H57").Value = 8380.933999999999
$ws.Range("J57").Value = 8380.933999999999
$ws.Range("L57").Value = 8380.933999999999
$ws.Range("N57").Value = -10020.934

$ws.Range("H82").Value = 0
$ws.Range("J82").Value = 0
$ws.Range("L82").Value = 0
$ws.Range("N82").Value = $null

$ws.Range("H85").Value = 0
$ws.Range("J85").Value = 0
$ws.Range("L85").Value = 0
$ws.Range("N85").Value = $null

$ws.Range("H132").Value = 7581940.5
$ws.Range("I132").Value = 15158708
$ws.Range("K132").Value = 45476124
$ws.Range("M132").Value = -45473594

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H55").Value = 140
$ws.Range("I55").Value = 165
$ws.Range("J55").Value = 90
$ws.Range("K55").Value = 165
$ws.Range("L55").Value = 90
$ws.Range("M55").Value = 8
$ws.Range("N55").Value = -436

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H2").Value = 39001
$ws.Range("J2").Value = 39001
$ws.Range("L2").Value = 39001
$ws.Range("N2").Value = -39225

$ws.Range("H5").Value = 0
$ws.Range("I5").Value = 0
$ws.Range("J5").Value = 0
$ws.Range("K5").Value = 0
$ws.Range("L5").Value = 0
$ws.Range("M5").Value = $null
$ws.Range("N5").Value = $null
